# IQC creator v1.0 20221228-003
#
# Row 3 (物料 10200054 / AAA-10200054-02) gets its IQC grading bumped from
# "B" to "A" for both version columns, and the two checkbox-style columns
# ("对勾1" / "对勾2") swap which one is marked (the "R" check-mark, rendered
# in Wingdings 2, moves from column N to column M; the empty box "□" moves
# the other way). The active sheet view is also updated to reflect where
# the author was working when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: grade columns J (IQC版本) and K (IQC_TB版本): "B" -> "A" ---
$ws.Range("J3").Value = "A"
$ws.Range("K3").Value = "A"

# --- Row 3: swap the checkbox cells M3 ("对勾1") and N3 ("对勾2") -------
# M3 currently holds the empty box "□" with no special formatting, while
# N3 holds the "R" check-glyph styled in Wingdings 2 (cell style index 5).
# After the edit, M3 should carry the check-glyph + style and N3 should
# carry the empty box with the default style. Routing the swap through a
# scratch cell via Range.Copy (rather than re-typing Value/Style) carries
# the exact existing cell style along with the cell content, so no new
# style entries are created.
$scratch = $ws.Range("BZ9999")
$ws.Range("M3").Copy($scratch)
$ws.Range("N3").Copy($ws.Range("M3"))
$scratch.Copy($ws.Range("N3"))
$scratch.ClearContents()

# --- Sheet view: scroll/selection reflects the author's last position ---
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollColumn = 6
    $excel.ActiveWindow.ScrollRow = 1
} catch {
    # Older/limited hosts may not expose window scrolling; selection below
    # is the important, persisted part of the view state.
}
$ws.Range("O7").Select()
